$d = $word.ActiveDocument
$r = $d.Content

# Simple forward sequential replace: finds the NEXT occurrence of $oldText
# starting from wherever the last Find left off, and replaces just that one
# occurrence (wdReplaceOne), leaving the run/formatting structure of any
# OTHER runs untouched.
function Replace-NextOccurrence($range, $oldText, $newText) {
    $found = $range.Find.Execute($oldText, $true, $true, $false, $false, $false, $true, 0, $false, $newText, 1)
    if (-not $found) {
        throw "Could not find next occurrence of '$oldText'"
    }
}

# Some table cells hold their value as TWO runs, e.g. "1.1 (3.4" + ")" -
# that second run exists separately (different w:rsidR) purely because of
# Word's edit history, but it shares the exact same formatting as the first
# run. A plain text replace on the first run causes this engine (like real
# Word often does) to coalesce the two equally-formatted adjacent runs into
# one - which would NOT match the target OOXML (the diff keeps both runs).
# To stop the auto-merge we momentarily make the first run's formatting
# differ from the second run (toggle Bold on, edit the text, toggle Bold
# back off) so the two runs are never "equal" at the moment of the edit.
function Set-FirstRunTextKeepingTrailingRun($table, $row, $col, $newFirstPart) {
    $cell = $table.Cell($row, $col)
    $cellRange = $cell.Range
    $charCount = $cellRange.Characters.Count   # includes the trailing cell mark
    $textLen = $charCount - 1                  # visible text length only
    $run1Start = $cellRange.Start
    $run2Start = $cellRange.Start + $textLen - 1
    $run2End = $cellRange.Start + $textLen

    $run1 = $d.Range($run1Start, $run2Start)
    $run1.Bold = 1
    $run1.Text = $newFirstPart

    $newRun1End = $run1Start + $newFirstPart.Length
    $run1Again = $d.Range($run1Start, $newRun1End)
    $run1Again.Bold = 0
}

Replace-NextOccurrence $r "03-Sep-2021, 14:36 " "08-Sep-2021, 15:21 "
Replace-NextOccurrence $r "TAC S1 M1" "TAC S1 D4"
Replace-NextOccurrence $r "40" "280"
Replace-NextOccurrence $r "24.1" "168.9"
Replace-NextOccurrence $r "2+2" "1+1"
Replace-NextOccurrence $r "5.3" "4.4"
Replace-NextOccurrence $r "80" "150"
Replace-NextOccurrence $r "276041.6" "1894184.5"
Replace-NextOccurrence $r "48" "304"
Replace-NextOccurrence $r "37.8" "38.3"
Replace-NextOccurrence $r "2" "1"
Replace-NextOccurrence $r "5" "6.5"
Replace-NextOccurrence $r "100" "200"
Replace-NextOccurrence $r "155081.6" "1047464.5"
Replace-NextOccurrence $r "277.2" "1872.5"
Replace-NextOccurrence $r "100" "150"
Replace-NextOccurrence $r "20" "40"
Replace-NextOccurrence $r "5.2" "9.1"

# "Absorbent pump rating" / "Refrigerant pump rating" rows: value cell is
# split into two runs ("1.1 (3.4" + ")"); preserve that split.
$tbl2 = $d.Tables.Item(2)
Set-FirstRunTextKeepingTrailingRun $tbl2 43 4 "3 (8"
Set-FirstRunTextKeepingTrailingRun $tbl2 44 4 "0.3 (1.4"

Replace-NextOccurrence $r "2800" "4120"
Replace-NextOccurrence $r "1450" "2130"
Replace-NextOccurrence $r "2250" "3045"
Replace-NextOccurrence $r "3.2" "9.8"
Replace-NextOccurrence $r "2.7" "7.5"
Replace-NextOccurrence $r "3" "9.4"
Replace-NextOccurrence $r "4.6" "15.2"
Replace-NextOccurrence $r "2700" "3640"

Write-Output "Done"
